$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "29.149.07";      E = "  +0.16%  " }
    @{ Row = 3;  D = "1.824.77";       E = "  -0.34%  " }
    @{ Row = 4;  D = "0.9988";         E = "  -0.02%  " }
    @{ Row = 5;  D = "241.57";         E = "  -0.64%  " }
    @{ Row = 6;  D = "0.6171";         E = "  -1.78%  " }
    @{ Row = 7;  D = "1.000";          E = "  -0.02%  " }
    @{ Row = 8;  D = "0.07342";        E = "  -2.39%  " }
    @{ Row = 9;  D = "0.2885";         E = "  -1.25%  " }
    @{ Row = 10; D = "22.98";          E = "  -0.91%  " }
    @{ Row = 11; D = "0.07668";        E = "  -0.09%  " }
    @{ Row = 12; D = "1.826.05";       E = "  -0.22%  " }
    @{ Row = 13; D = "4.950";          E = "  -1.12%  " }
    @{ Row = 14; D = "0.6610";         E = "  -1.04%  " }
    @{ Row = 15; D = "81.89";          E = "  -1.07%  " }
    @{ Row = 16; D = "0.000008916";    E = "  -5.10%  " }
    @{ Row = 17; D = "5.838";          E = "  -2.32%  " }
    @{ Row = 18; D = "29.108.62";      E = "  +0.12%  " }
    @{ Row = 19; D = "2.070.17";       E = "  -0.31%  " }
    @{ Row = 20; D = "236.78";         E = "  +6.12%  " }
    @{ Row = 21; D = "12.42";          E = "  -1.11%  " }
    @{ Row = 22; D = "0.9999";         E = "  -0.20%  " }
    @{ Row = 23; D = "7.128";          E = "  +0.25%  " }
    @{ Row = 24; D = $null;            E = "  -0.03%  " }
    @{ Row = 25; D = "157.80";         E = "  -1.24%  " }
    @{ Row = 26; D = $null;            E = "  +1.40%  " }
    @{ Row = 27; D = "8.430";          E = "  -0.66%  " }
    @{ Row = 28; D = "17.60";          E = "  -1.50%  " }
    @{ Row = 29; D = "1.484";          E = "  -0.72%  " }
    @{ Row = 30; D = "0.05548";        E = "  -4.02%  " }
    @{ Row = 31; D = "4.088";          E = "  -0.07%  " }
    @{ Row = 32; D = "4.093";          E = "  -1.35%  " }
    @{ Row = 33; D = "1.203";          E = "  -0.40%  " }
    @{ Row = 34; D = "1.825";          E = "  -0.32%  " }
    @{ Row = 35; D = "0.7346";         E = "  -0.59%  " }
    @{ Row = 36; D = "1.132";          E = "  -0.39%  " }
    @{ Row = 37; D = "2.607";          E = "  -2.31%  " }
    @{ Row = 38; D = "2.835";          E = "  +2.48%  " }
    @{ Row = 39; D = "1.209.22";       E = "  -0.68%  " }
    @{ Row = 40; D = "0.01757";        E = "  -1.10%  " }
    @{ Row = 41; D = $null;            E = "  -2.48%  " }
    @{ Row = 42; D = "0.9010";         E = "  +1.39%  " }
    @{ Row = 43; D = "1.000";          E = "  -0.10%  " }
    @{ Row = 44; D = "101.40";         E = "  -0.67%  " }
    @{ Row = 45; D = "1.973.98";       E = "  -0.19%  " }
    @{ Row = 46; D = "0.00000000124"; E = "  +0.16%  " }
    @{ Row = 47; D = "64.58";          E = "  -1.53%  " }
    @{ Row = 48; D = "0.5077";         E = "  -0.32%  " }
    @{ Row = 49; D = "0.4004";         E = "  -1.40%  " }
    @{ Row = 50; D = "9.010";          E = "  +0.30%  " }
    @{ Row = 51; D = "0.05749";        E = "  -1.24%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    $dVal = $u.D

    if ($null -ne $dVal) {
        # Determine whether the textual price would be auto-parsed into a
        # number by Excel. If so, force the cell to Text format first so
        # that the stored value stays an exact string match (same as the
        # original inline string), otherwise just write the value.
        $isNumeric = $dVal -match '^[0-9]+(\.[0-9]+)?$'

        $dCell = $ws.Cells.Item($row, 4)
        if ($isNumeric) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $dVal
    }

    $ws.Cells.Item($row, 5).Value = $u.E
}
